# Update scripts with new tpm
# The NATMI pipeline was re-run with updated TPM data: a new sending cluster
# "ECs" now appears alongside the existing "FAPs" cluster, so the Ccl20-Cxcr3
# LR-pair sheet gains two more rows (ECs -> FAPs, ECs -> Resolving-Mac) and the
# derived-specificity columns for the pre-existing FAPs rows are recomputed to
# account for the extra cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs -> FAPs (was FAPs -> FAPs) ---------------------------------
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl20"
$ws.Range("C2").Value = "Cxcr3"
$ws.Range("D2").Value = "FAPs"

$ws.Range("G2").Value = 0.07886166666666666
$ws.Range("H2").Value = 0.236585
$ws.Range("I2").Value = 0.5304501264551309
$ws.Range("J2").Value = 0.5304501264551309
$ws.Range("O2").Value = 0.02773017886769741
$ws.Range("P2").Value = 0.02773017886769741
$ws.Range("Q2").Value = 0.004096180115555555
$ws.Range("R2").Value = 0.03686562104
$ws.Range("S2").Value = 0.01470947688699349
$ws.Range("T2").Value = 0.01470947688699349

# --- Row 3: ECs -> Resolving-Mac (was FAPs -> Resolving-Mac) --------------
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ccl20"
$ws.Range("C3").Value = "Cxcr3"
$ws.Range("D3").Value = "Resolving-Mac"

$ws.Range("G3").Value = 0.07886166666666666
$ws.Range("H3").Value = 0.236585
$ws.Range("I3").Value = 0.5304501264551309
$ws.Range("J3").Value = 0.5304501264551309
$ws.Range("M3").Value = 1.821156333333333
$ws.Range("N3").Value = 5.463469
$ws.Range("O3").Value = 0.9722698211323025
$ws.Range("P3").Value = 0.9722698211323026
$ws.Range("Q3").Value = 0.1436194237072222
$ws.Range("R3").Value = 1.292574813365
$ws.Range("S3").Value = 0.5157406495681374
$ws.Range("T3").Value = 0.5157406495681375

# --- Row 4 (new): FAPs -> FAPs ---------------------------------------------
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ccl20"
$ws.Range("C4").Value = "Cxcr3"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.06980766666666667
$ws.Range("H4").Value = 0.209423
$ws.Range("I4").Value = 0.4695498735448692
$ws.Range("J4").Value = 0.4695498735448692
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05194133333333333
$ws.Range("N4").Value = 0.155824
$ws.Range("O4").Value = 0.02773017886769741
$ws.Range("P4").Value = 0.02773017886769741
$ws.Range("Q4").Value = 0.003625903283555556
$ws.Range("R4").Value = 0.032633129552
$ws.Range("S4").Value = 0.01302070198070392
$ws.Range("T4").Value = 0.01302070198070392

# --- Row 5 (new): FAPs -> Resolving-Mac ------------------------------------
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ccl20"
$ws.Range("C5").Value = "Cxcr3"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.06980766666666667
$ws.Range("H5").Value = 0.209423
$ws.Range("I5").Value = 0.4695498735448692
$ws.Range("J5").Value = 0.4695498735448692
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.821156333333333
$ws.Range("N5").Value = 5.463469
$ws.Range("O5").Value = 0.9722698211323025
$ws.Range("P5").Value = 0.9722698211323026
$ws.Range("Q5").Value = 0.1271306742652222
$ws.Range("R5").Value = 1.144176068387
$ws.Range("S5").Value = 0.4565291715641653
$ws.Range("T5").Value = 0.4565291715641653
